# Add a "Units" column to the "Concentrations" sheet, reflecting the new
# `units` attribute added to wc_lang.Concentration.
#
# Before:  Species | Value | Comments  | References
# After:   Species | Value | Units     | Comments | References
#
# Concentration rows expressed per liquid volume (the specie_*[e]/[c] rows)
# are recorded in molar units ("M"); the two H2O rows are recorded as
# particle counts ("molecules").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concentrations")

# Insert a new blank column in front of the existing "Comments" column (C),
# pushing Comments -> D and References -> E.
$ws.Columns.Item(3).Insert()

$ws.Cells.Item(1, 3).Value = "Units"

$ws.Cells.Item(2, 3).Value = "M"
$ws.Cells.Item(3, 3).Value = "M"
$ws.Cells.Item(4, 3).Value = "M"
$ws.Cells.Item(5, 3).Value = "M"
$ws.Cells.Item(6, 3).Value = "M"
$ws.Cells.Item(7, 3).Value = "M"
$ws.Cells.Item(8, 3).Value = "molecules"
$ws.Cells.Item(9, 3).Value = "molecules"

# Rows 8-9 (the H2O rows) use a different style than column B in those rows,
# so copy the number/text format used by the rest of the new Units column
# (as seen in C2) down into C8:C9.
$ws.Range("C2").Copy()
$ws.Range("C8:C9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Make Concentrations the active sheet/tab and leave the selection on C8,
# matching the state the workbook was saved in.
$ws.Activate()
$ws.Range("C8").Select()
